# Append a new list item after the paragraph that ends with
# "...because it is not mandatory and deadline issues." and move the
# "_GoBack" bookmark (if present) to the end of the new paragraph.

$d = $word.ActiveDocument

# --- Locate the paragraph we need to split after -----------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*mandatory and deadline issues.*") {
        $targetPara = $p
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the target paragraph (...deadline issues.)"
}

$newText = "Made the ticket price button a bit more obvious."

# Position right after the last real character of the paragraph, i.e.
# right before its paragraph mark (and before any bookmark sitting there).
$paraEnd = $targetPara.Range.End - 1

# A temporary marker is appended after the new text so that the position
# where the bookmark needs to be created is never the very last
# character(s) of the document content while we create it (a collapsed
# bookmark placed exactly at Content.End misbehaves), we trim it right
# after placing the bookmark.
$sentinel = "~~~TMP_SENTINEL~~~"

$insertPoint = $d.Range($paraEnd, $paraEnd)
$insertPoint.InsertAfter([char]13 + $newText + $sentinel)

# Position right after the newly inserted text (= where the bookmark
# used to be relative to the old paragraph end).
$bmPos = $paraEnd + 1 + $newText.Length
$bmRange = $d.Range($bmPos, $bmPos)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary sentinel text again.
$trimRange = $d.Range($bmPos, $bmPos + $sentinel.Length)
$trimRange.Delete()
